# Balance Fix (보우미터 사거리) - adjust UseRange values and page setup

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update UseRange (column E) values for the skill table
$ws.Range("E4").Value = 10
$ws.Range("E5").Value = 15
$ws.Range("E6").Value = 15

# Configure page setup (paper size A4, portrait orientation)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Move the active selection to J5
[void]$ws.Range("J5").Select()
